$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
if (-not $ws) { $ws = $wb.ActiveSheet }

# New time-tracking entry in row 34 (2020-08-17, 2.5h)
$ws.Range("A34").Value = 44060
$ws.Range("B34").Value = 2.5
$ws.Range("C34").Value = "usernamen generoimista ja cookieiden käyttöä localstoragen sijasta."

# The description wraps onto two lines at this column width, so the row grows
$ws.Rows.Item(34).RowHeight = 30

# Scroll the view down to the new row and select the newly filled cell,
# matching where the author's cursor ended up after the edit
$win = $excel.ActiveWindow
$win.ScrollRow = 28
$win.ScrollColumn = 1
$win.Left = -27360
$win.Top = 720
$ws.Range("C34").Select()
